$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "56.578.42"
$ws.Range("E2").Value = "  +4.38%  "
Set-TextValue $ws "D3" "3.009.78"
$ws.Range("E3").Value = "  +5.07%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws "D5" "508.70"
$ws.Range("E5").Value = "  +8.13%  "
Set-TextValue $ws "D6" "136.77"
$ws.Range("E6").Value = "  +8.88%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +8.10%  "
Set-TextValue $ws "D9" "7.63"
$ws.Range("E9").Value = "  +15.61%  "
$ws.Range("E10").Value = "  +14.14%  "
$ws.Range("E11").Value = "  +7.75%  "
$ws.Range("E12").Value = "  +4.67%  "
Set-TextValue $ws "D13" "3.524.94"
$ws.Range("E13").Value = "  +4.96%  "
Set-TextValue $ws "D14" "25.74"
$ws.Range("E14").Value = "  +11.01%  "
$ws.Range("E15").Value = "  +15.52%  "
Set-TextValue $ws "D16" "56.614.70"
$ws.Range("E16").Value = "  +4.24%  "
Set-TextValue $ws "D17" "3.009.81"
$ws.Range("E17").Value = "  +5.02%  "
$ws.Range("E18").Value = "  +9.92%  "
Set-TextValue $ws "D19" "12.52"
$ws.Range("E19").Value = "  +9.73%  "
Set-TextValue $ws "D20" "7.87"
$ws.Range("E20").Value = "  +12.20%  "
Set-TextValue $ws "D21" "327.79"
$ws.Range("E21").Value = "  +11.75%  "
Set-TextValue $ws "D22" "0.999"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +8.04%  "
Set-TextValue $ws "D24" "62.52"
$ws.Range("E24").Value = "  +6.67%  "
$ws.Range("E25").Value = "  +9.39%  "
Set-TextValue $ws "D26" "1.00"
$ws.Range("E26").Value = "  +0.92%  "
Set-TextValue $ws "D27" "0.0₃0919"
$ws.Range("E27").Value = "  +14.08%  "
$ws.Range("E28").Value = "  +6.53%  "
$ws.Range("E29").Value = "  +12.73%  "
Set-TextValue $ws "D30" "1.25"
$ws.Range("E30").Value = "  +11.12%  "
$ws.Range("E31").Value = "  +10.01%  "
Set-TextValue $ws "D32" "20.63"
$ws.Range("E32").Value = "  +9.84%  "
Set-TextValue $ws "D33" "156.13"
$ws.Range("E33").Value = "  +16.19%  "
$ws.Range("E34").Value = "  +7.64%  "
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("E36").Value = "  +4.83%  "
Set-TextValue $ws "D37" "0.0675"
$ws.Range("E37").Value = "  +9.46%  "
Set-TextValue $ws "D38" "23.72"
$ws.Range("E38").Value = "  +3.20%  "
Set-TextValue $ws "D39" "3.044.80"
$ws.Range("E39").Value = "  +5.25%  "
Set-TextValue $ws "D40" "36.60"
$ws.Range("E40").Value = "  +4.71%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +8.11%  "
Set-TextValue $ws "D43" "2.267.07"
$ws.Range("E43").Value = "  +10.98%  "
$ws.Range("E44").Value = "  +6.03%  "
$ws.Range("E45").Value = "  +7.66%  "
$ws.Range("E46").Value = "  +6.93%  "
$ws.Range("E47").Value = "  +24.73%  "
Set-TextValue $ws "D48" "0.0236"
$ws.Range("E48").Value = "  +11.57%  "
$ws.Range("E49").Value = "  +8.65%  "
Set-TextValue $ws "D50" "19.17"
$ws.Range("E50").Value = "  +7.25%  "
$ws.Range("E51").Value = "  +11.15%  "
